# =====================================================================
# Edit script: add "2022-Q3" sheet with new fund data, and update the
# "总计" (totals) summary sheet with a new leading row for 2022-Q3.
# =====================================================================

$wb = $excel.ActiveWorkbook
$wsTotal = $wb.Worksheets.Item(1)

# ---- Build new sheet "2022-Q3" ----
$ws1 = $wb.Worksheets.Item(1)
$wsQ3 = $wb.Worksheets.Add($null, $ws1)
$wsQ3.Name = "2022-Q3"

# Look up the (still-present, not-yet-renamed) "2022-Q2" sheet AFTER the
# structural Add()/Name= changes above -- object/name lookups taken before
# a sheet insertion can resolve to a stale reference.
$wsOldQ2 = $wb.Worksheets.Item("2022-Q2")

# Protect numeric-looking text columns (B, D, E, F, G for rows 2-28) from auto-number conversion
$wsQ3.Range("B2:B28").NumberFormat = "@"
$wsQ3.Range("D2:G28").NumberFormat = "@"

# Row 1
$wsQ3.Range("B1").Value = "基金代码"
$wsQ3.Range("C1").Value = "基金名称"
$wsQ3.Range("D1").Value = "基金规模"
$wsQ3.Range("E1").Value = "股票总仓位"
$wsQ3.Range("F1").Value = "仓位占比"
$wsQ3.Range("G1").Value = "持有市值(亿元)"
$wsQ3.Range("H1").Value = "仓位排名"
# Row 2
$wsQ3.Range("A2").Value = 0
$wsQ3.Range("B2").Value = "159941"
$wsQ3.Range("C2").Value = "广发纳斯达克100ETF（QDII）"
$wsQ3.Range("D2").Value = "106.15"
$wsQ3.Range("E2").Value = "91.14"
$wsQ3.Range("F2").Value = "2.31"
$wsQ3.Range("G2").Value = "2.4521"
$wsQ3.Range("H2").Value = 8
# Row 3
$wsQ3.Range("A3").Value = 1
$wsQ3.Range("B3").Value = "000906"
$wsQ3.Range("C3").Value = "广发全球精选股票（QDII）美元现汇"
$wsQ3.Range("D3").Value = "21.88"
$wsQ3.Range("E3").Value = "79.27"
$wsQ3.Range("F3").Value = "6.99"
$wsQ3.Range("G3").Value = "1.5294"
$wsQ3.Range("H3").Value = 2
# Row 4
$wsQ3.Range("A4").Value = 2
$wsQ3.Range("B4").Value = "270023"
$wsQ3.Range("C4").Value = "广发全球精选股票（QDII）"
$wsQ3.Range("D4").Value = "21.88"
$wsQ3.Range("E4").Value = "79.27"
$wsQ3.Range("F4").Value = "6.99"
$wsQ3.Range("G4").Value = "1.5294"
$wsQ3.Range("H4").Value = 2
# Row 5
$wsQ3.Range("A5").Value = 3
$wsQ3.Range("B5").Value = "011421"
$wsQ3.Range("C5").Value = "广发全球科技三个月定期开放混合（QDII）美元 A"
$wsQ3.Range("D5").Value = "22.73"
$wsQ3.Range("E5").Value = "90.00"
$wsQ3.Range("F5").Value = "5.38"
$wsQ3.Range("G5").Value = "1.2229"
$wsQ3.Range("H5").Value = 4
# Row 6
$wsQ3.Range("A6").Value = 4
$wsQ3.Range("B6").Value = "011420"
$wsQ3.Range("C6").Value = "广发全球科技三个月定期开放混合（QDII）人民币 A"
$wsQ3.Range("D6").Value = "22.73"
$wsQ3.Range("E6").Value = "90.00"
$wsQ3.Range("F6").Value = "5.38"
$wsQ3.Range("G6").Value = "1.2229"
$wsQ3.Range("H6").Value = 4
# Row 7
$wsQ3.Range("A7").Value = 5
$wsQ3.Range("B7").Value = "513100"
$wsQ3.Range("C7").Value = "国泰纳斯达克100（QDII-ETF）"
$wsQ3.Range("D7").Value = "46.54"
$wsQ3.Range("E7").Value = "91.35"
$wsQ3.Range("F7").Value = "2.32"
$wsQ3.Range("G7").Value = "1.0797"
$wsQ3.Range("H7").Value = 8
# Row 8
$wsQ3.Range("A8").Value = 6
$wsQ3.Range("B8").Value = "040047"
$wsQ3.Range("C8").Value = "华安纳斯达克100指数（QDII）美元现钞A"
$wsQ3.Range("D8").Value = "24.52"
$wsQ3.Range("E8").Value = "92.09"
$wsQ3.Range("F8").Value = "2.34"
$wsQ3.Range("G8").Value = "0.5738"
$wsQ3.Range("H8").Value = 8
# Row 9
$wsQ3.Range("A9").Value = 7
$wsQ3.Range("B9").Value = "040048"
$wsQ3.Range("C9").Value = "华安纳斯达克100指数（QDII）美元现汇A"
$wsQ3.Range("D9").Value = "24.52"
$wsQ3.Range("E9").Value = "92.09"
$wsQ3.Range("F9").Value = "2.34"
$wsQ3.Range("G9").Value = "0.5738"
$wsQ3.Range("H9").Value = 8
# Row 10
$wsQ3.Range("A10").Value = 8
$wsQ3.Range("B10").Value = "040046"
$wsQ3.Range("C10").Value = "华安纳斯达克100指数（QDII）人民币A"
$wsQ3.Range("D10").Value = "22.21"
$wsQ3.Range("E10").Value = "92.09"
$wsQ3.Range("F10").Value = "2.34"
$wsQ3.Range("G10").Value = "0.5197"
$wsQ3.Range("H10").Value = 8
# Row 11
$wsQ3.Range("A11").Value = 9
$wsQ3.Range("B11").Value = "160213"
$wsQ3.Range("C11").Value = "国泰纳斯达克100指数（QDII）"
$wsQ3.Range("D11").Value = "15.14"
$wsQ3.Range("E11").Value = "85.81"
$wsQ3.Range("F11").Value = "2.22"
$wsQ3.Range("G11").Value = "0.3361"
$wsQ3.Range("H11").Value = 8
# Row 12
$wsQ3.Range("A12").Value = 10
$wsQ3.Range("B12").Value = "000834"
$wsQ3.Range("C12").Value = "大成纳斯达克100指数（QDII）"
$wsQ3.Range("D12").Value = "14.15"
$wsQ3.Range("E12").Value = "85.22"
$wsQ3.Range("F12").Value = "2.17"
$wsQ3.Range("G12").Value = "0.3071"
$wsQ3.Range("H12").Value = 8
# Row 13
$wsQ3.Range("A13").Value = 11
$wsQ3.Range("B13").Value = "011423"
$wsQ3.Range("C13").Value = "广发全球科技三个月定期开放混合（QDII）美元 C"
$wsQ3.Range("D13").Value = "5.27"
$wsQ3.Range("E13").Value = "90.00"
$wsQ3.Range("F13").Value = "5.38"
$wsQ3.Range("G13").Value = "0.2835"
$wsQ3.Range("H13").Value = 4
# Row 14
$wsQ3.Range("A14").Value = 12
$wsQ3.Range("B14").Value = "011422"
$wsQ3.Range("C14").Value = "广发全球科技三个月定期开放混合（QDII）人民币 C"
$wsQ3.Range("D14").Value = "5.27"
$wsQ3.Range("E14").Value = "90.00"
$wsQ3.Range("F14").Value = "5.38"
$wsQ3.Range("G14").Value = "0.2835"
$wsQ3.Range("H14").Value = 4
# Row 15
$wsQ3.Range("A15").Value = 13
$wsQ3.Range("B15").Value = "513300"
$wsQ3.Range("C15").Value = "华夏纳斯达克100ETF（QDII）"
$wsQ3.Range("D15").Value = "11.08"
$wsQ3.Range("E15").Value = "97.32"
$wsQ3.Range("F15").Value = "2.48"
$wsQ3.Range("G15").Value = "0.2748"
$wsQ3.Range("H15").Value = 2
# Row 16
$wsQ3.Range("A16").Value = 14
$wsQ3.Range("B16").Value = "000043"
$wsQ3.Range("C16").Value = "嘉实美国成长股票（QDII）人民币"
$wsQ3.Range("D16").Value = "12.41"
$wsQ3.Range("E16").Value = "92.80"
$wsQ3.Range("F16").Value = "1.77"
$wsQ3.Range("G16").Value = "0.2197"
$wsQ3.Range("H16").Value = 7
# Row 17
$wsQ3.Range("A17").Value = 15
$wsQ3.Range("B17").Value = "000044"
$wsQ3.Range("C17").Value = "嘉实美国成长股票（QDII）美元现汇"
$wsQ3.Range("D17").Value = "12.41"
$wsQ3.Range("E17").Value = "92.80"
$wsQ3.Range("F17").Value = "1.77"
$wsQ3.Range("G17").Value = "0.2197"
$wsQ3.Range("H17").Value = 7
# Row 18
$wsQ3.Range("A18").Value = 16
$wsQ3.Range("B18").Value = "161128"
$wsQ3.Range("C18").Value = "易方达标普信息科技指数（QDII-LOF）人民币"
$wsQ3.Range("D18").Value = "4.99"
$wsQ3.Range("E18").Value = "91.96"
$wsQ3.Range("F18").Value = "3.52"
$wsQ3.Range("G18").Value = "0.1756"
$wsQ3.Range("H18").Value = 3
# Row 19
$wsQ3.Range("A19").Value = 17
$wsQ3.Range("B19").Value = "012868"
$wsQ3.Range("C19").Value = "易方达标普信息科技指数（QDII-LOF）人民币 C"
$wsQ3.Range("D19").Value = "4.99"
$wsQ3.Range("E19").Value = "91.96"
$wsQ3.Range("F19").Value = "3.52"
$wsQ3.Range("G19").Value = "0.1756"
$wsQ3.Range("H19").Value = 3
# Row 20
$wsQ3.Range("A20").Value = 18
$wsQ3.Range("B20").Value = "003722"
$wsQ3.Range("C20").Value = "易方达纳斯达克100指数美元（QDII-LOF）A"
$wsQ3.Range("D20").Value = "7.72"
$wsQ3.Range("E20").Value = "90.67"
$wsQ3.Range("F20").Value = "2.26"
$wsQ3.Range("G20").Value = "0.1745"
$wsQ3.Range("H20").Value = 8
# Row 21
$wsQ3.Range("A21").Value = 19
$wsQ3.Range("B21").Value = "161130"
$wsQ3.Range("C21").Value = "易方达纳斯达克100指数人民币（QDII-LOF）"
$wsQ3.Range("D21").Value = "7.72"
$wsQ3.Range("E21").Value = "90.67"
$wsQ3.Range("F21").Value = "2.26"
$wsQ3.Range("G21").Value = "0.1745"
$wsQ3.Range("H21").Value = 8
# Row 22
$wsQ3.Range("A22").Value = 20
$wsQ3.Range("B22").Value = "003721"
$wsQ3.Range("C22").Value = "易方达标普信息科技指数（QDII-LOF）美元A"
$wsQ3.Range("D22").Value = "4.84"
$wsQ3.Range("E22").Value = "91.96"
$wsQ3.Range("F22").Value = "3.52"
$wsQ3.Range("G22").Value = "0.1704"
$wsQ3.Range("H22").Value = 3
# Row 23
$wsQ3.Range("A23").Value = 21
$wsQ3.Range("B23").Value = "014978"
$wsQ3.Range("C23").Value = "华安纳斯达克100指数（QDII）人民币C"
$wsQ3.Range("D23").Value = "2.31"
$wsQ3.Range("E23").Value = "92.09"
$wsQ3.Range("F23").Value = "2.34"
$wsQ3.Range("G23").Value = "0.0541"
$wsQ3.Range("H23").Value = 8
# Row 24
$wsQ3.Range("A24").Value = 22
$wsQ3.Range("B24").Value = "159632"
$wsQ3.Range("C24").Value = "华安纳斯达克100ETF（QDII）"
$wsQ3.Range("D24").Value = "1.51"
$wsQ3.Range("E24").Value = "89.05"
$wsQ3.Range("F24").Value = "2.31"
$wsQ3.Range("G24").Value = "0.0349"
$wsQ3.Range("H24").Value = 8
# Row 25
$wsQ3.Range("A25").Value = 23
$wsQ3.Range("B25").Value = "005698"
$wsQ3.Range("C25").Value = "华夏全球科技先锋混合（QDII）"
$wsQ3.Range("D25").Value = "0.59"
$wsQ3.Range("E25").Value = "86.79"
$wsQ3.Range("F25").Value = "4.67"
$wsQ3.Range("G25").Value = "0.0276"
$wsQ3.Range("H25").Value = 7
# Row 26
$wsQ3.Range("A26").Value = 24
$wsQ3.Range("B26").Value = "012869"
$wsQ3.Range("C26").Value = "易方达标普信息科技指数（QDII-LOF）美元 C"
$wsQ3.Range("D26").Value = "0.15"
$wsQ3.Range("E26").Value = "91.96"
$wsQ3.Range("F26").Value = "3.52"
$wsQ3.Range("G26").Value = "0.0053"
$wsQ3.Range("H26").Value = 3
# Row 27
$wsQ3.Range("A27").Value = 25
$wsQ3.Range("B27").Value = "012871"
$wsQ3.Range("C27").Value = "易方达纳斯达克100指数美元（QDII-LOF）C"
$wsQ3.Range("D27").Value = "0.18"
$wsQ3.Range("E27").Value = "90.67"
$wsQ3.Range("F27").Value = "2.26"
$wsQ3.Range("G27").Value = "0.0041"
$wsQ3.Range("H27").Value = 8
# Row 28
$wsQ3.Range("A28").Value = 26
$wsQ3.Range("B28").Value = "012870"
$wsQ3.Range("C28").Value = "易方达纳斯达克100指数人民币（QDII-LOF）C"
$wsQ3.Range("D28").Value = "0.18"
$wsQ3.Range("E28").Value = "90.67"
$wsQ3.Range("F28").Value = "2.26"
$wsQ3.Range("G28").Value = "0.0041"
$wsQ3.Range("H28").Value = 8

# Reset number format for the text columns back to General (keeps stored type as text, clears the dangling format)
$wsQ3.Range("B2:B28").Style = "Normal"
$wsQ3.Range("D2:G28").Style = "Normal"


# ---- Apply header style (bold + border, matches existing quarter sheets) ----
$wsOldQ2.Range("B1:H1").Copy()
$wsQ3.Range("B1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Apply row-index column (A) style (bold + border), matches existing quarter sheets ----
$wsOldQ2.Range("A2").Copy()
$wsQ3.Range("A2:A28").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Update "总计" sheet: insert new leading row for 2022-Q3, shift the rest down ----
$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 27
$wsTotal.Range("D2").Value = 13.63

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q2"
$wsTotal.Range("C3").Value = 36
$wsTotal.Range("D3").Value = 15.9

$wsTotal.Range("A4").Value = 2
$wsTotal.Range("B4").Value = "2022-Q1"
$wsTotal.Range("C4").Value = 41
$wsTotal.Range("D4").Value = 38.4

$wsTotal.Range("A5").Value = 3
$wsTotal.Range("B5").Value = "2021-Q4"
$wsTotal.Range("C5").Value = 37
$wsTotal.Range("D5").Value = 43.76

$wsTotal.Range("A6").Value = 4
$wsTotal.Range("B6").Value = "2021-Q1"
$wsTotal.Range("C6").Value = 24
$wsTotal.Range("D6").Value = 11.88

# Row 6 is brand new territory (previous used range stopped at row 5) so it
# needs to explicitly inherit the row-index column style from row 5.
$wsTotal.Range("A5").Copy()
$wsTotal.Range("A6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

